# "Use raw SayCan prompt" — fill in the raw-prompt SayCan results (and the
# other already-started metrics) for trials 1 and 2 in columns B:Y / B:S.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Trial 1): B2:Y2, all 1 except T2 ("0 - SayCan") which is 0
$ws.Range("B2:S2").Value = 1
$ws.Range("T2").Value = 0
$ws.Range("U2:Y2").Value = 1

# Row 3 (Trial 2): B3:S3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0

# Update the view: scroll so column K is the left-most visible column,
# and leave the active selection on T3.
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 1
$ws.Range("T3").Select()
